$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 490, pushing existing rows 490-529 down to 491-530
$ws.Rows.Item(490).Insert()

# Populate the newly inserted row 490 with the new record's data
$ws.Cells.Item(490, 1).Value = 3
$ws.Cells.Item(490, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(490, 3).Value = "Coquimbo"
$ws.Cells.Item(490, 4).Value = 45013
$ws.Cells.Item(490, 4).NumberFormat = $ws.Cells.Item(491, 4).NumberFormat()
$ws.Cells.Item(490, 5).Value = 5
$ws.Cells.Item(490, 6).Value = 100112043
$ws.Cells.Item(490, 7).Value = "Pepino ensalada"
$ws.Cells.Item(490, 8).Value = "Sin especificar"
$ws.Cells.Item(490, 9).Value = "Primera"
$ws.Cells.Item(490, 10).Value = 130
$ws.Cells.Item(490, 11).Value = 10000
$ws.Cells.Item(490, 12).Value = 11000
$ws.Cells.Item(490, 13).Value = 10462
$ws.Cells.Item(490, 14).Value = "`$/caja 60 unidades"
$ws.Cells.Item(490, 15).Value = "Limache"
$ws.Cells.Item(490, 16).Value = 174
$ws.Cells.Item(490, 17).Value = 60
$ws.Cells.Item(490, 18).Value = "Hortaliza"
